$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.476.23'
$ws.Range('E2').Value = '  +1.84%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.855.07'
$ws.Range('E3').Value = '  +1.15%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9995'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.06'
$ws.Range('E5').Value = '  -0.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6938'
$ws.Range('E6').Value = '  +0.60%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3062'
$ws.Range('E8').Value = '  +0.33%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07658'
$ws.Range('E9').Value = '  -0.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.58'
$ws.Range('E10').Value = '  +0.17%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07772'
$ws.Range('E11').Value = '  -0.34%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.138'
$ws.Range('E12').Value = '  +1.16%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.852.55'
$ws.Range('E13').Value = '  +0.63%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6908'
$ws.Range('E14').Value = '  +1.59%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '90.75'
$ws.Range('E15').Value = '  +0.28%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.365'
$ws.Range('E16').Value = '  -1.10%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '29.466.77'
$ws.Range('E17').Value = '  +1.76%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008267'
$ws.Range('E18').Value = '  -0.86%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.101.98'
$ws.Range('E19').Value = '  +0.84%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '237.51'
$ws.Range('E20').Value = '  -2.40%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.72'
$ws.Range('E21').Value = '  +0.06%  '
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.680'
$ws.Range('E23').Value = '  +2.77%  '
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('E25').Value = '  +1.57%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.896'
$ws.Range('E26').Value = '  +1.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '159.64'
$ws.Range('E27').Value = '  -1.58%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.24'
$ws.Range('E28').Value = '  +0.20%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.538'
$ws.Range('E29').Value = '  -1.09%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.243'
$ws.Range('E30').Value = '  +0.75%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.153'
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.196'
$ws.Range('E32').Value = '  +1.60%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05101'
$ws.Range('E33').Value = '  -0.43%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7704'
$ws.Range('E34').Value = '  +0.55%  '
$ws.Range('E35').Value = '  +2.35%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.149'
$ws.Range('E36').Value = '  +0.47%  '
$ws.Range('E37').Value = '  +0.09%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.335.34'
$ws.Range('E38').Value = '  +8.18%  '
$ws.Range('E39').Value = '  +1.04%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.714'
$ws.Range('E40').Value = '  +0.68%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9620'
$ws.Range('E41').Value = '  +4.12%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '106.23'
$ws.Range('E42').Value = '  -2.16%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.844'
$ws.Range('E43').Value = '  +0.29%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.000'
$ws.Range('E44').Value = '  +0.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '9.805'
$ws.Range('E45').Value = '  +2.25%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000125'
$ws.Range('E46').Value = '  +1.84%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.000.18'
$ws.Range('E47').Value = '  +0.87%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5221'
$ws.Range('E48').Value = '  +0.97%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.781'
$ws.Range('E49').Value = '  +1.98%  '
$ws.Range('E51').Value = '  +0.63%  '
